$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.301.78"
$ws.Range("E2").Value = "  +8.64%  "
$ws.Range("D3").Value = "2.585.26"
$ws.Range("E3").Value = "  +6.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +15.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "578.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.203"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +25.04%  "
$ws.Range("D10").Value = "2.587.51"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  +8.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.77"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +9.89%  "
$ws.Range("D15").Value = "74.145.87"
$ws.Range("E15").Value = "  +8.57%  "
$ws.Range("D16").Value = "3.046.85"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +13.19%  "
$ws.Range("D18").Value = "2.573.34"
$ws.Range("E18").Value = "  +6.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +25.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +11.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +20.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.31%  "
$ws.Range("E26").Value = "  +12.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +11.15%  "
$ws.Range("D28").Value = "2.716.60"
$ws.Range("E28").Value = "  +6.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  +14.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "499.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +17.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +17.48%  "
$ws.Range("E34").Value = "  +5.93%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +12.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.78"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.90"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.67"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +12.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.320"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +19.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "38.86"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.49"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +12.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0828"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +16.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.518"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0967"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.82%  "
